$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115; existing rows 115-143 shift down to 116-144.
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with the new weekly data point.
$ws.Cells.Item(115, 1).Value = 1
$ws.Cells.Item(115, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(115, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(115, 4).Value = 44876
$ws.Cells.Item(115, 5).Value = 15
$ws.Cells.Item(115, 6).Value = 100112042
$ws.Cells.Item(115, 7).Value = "Locoto"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 130
$ws.Cells.Item(115, 11).Value = 21000
$ws.Cells.Item(115, 12).Value = 22000
$ws.Cells.Item(115, 13).Value = 21500
$ws.Cells.Item(115, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(115, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(115, 16).Value = 1075
$ws.Cells.Item(115, 17).Value = 20
$ws.Cells.Item(115, 18).Value = "Hortaliza"
